$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.147.35"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "1.679.15"
$ws.Range("E3").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.63%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'22.77"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +6.14%  "
$ws.Range("D9").Value = "'0.260"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.37%  "
$ws.Range("D10").Value = "'0.0621"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "'0.0891"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").Value = "1.916.62"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.678.62"
$ws.Range("E13").Value = "  -0.46%  "
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").Value = "'66.55"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").Value = "27.127.23"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'234.72"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").Value = "'7.88"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("E20").Value = "  +0.36%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  +1.43%  "
$ws.Range("E23").Value = "  +2.69%  "
$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("D25").Value = "'149.07"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "'7.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.41%  "
$ws.Range("D27").Value = "'16.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("E28").Value = "  -0.12%  "
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  +0.39%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "1.541.96"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("E34").Value = "  +0.63%  "
$ws.Range("E35").Value = "  -4.02%  "
$ws.Range("D36").Value = "'0.607"
$ws.Range("D36").ClearFormats()
$ws.Range("E37").Value = "  +2.69%  "
$ws.Range("E40").Value = "  +2.33%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'69.53"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.41%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.78"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "1.824.89"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "'89.84"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.67%  "
$ws.Range("E48").Value = "  +6.58%  "
$ws.Range("E49").Value = "  +2.61%  "
$ws.Range("D50").Value = "'8.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.76%  "
$ws.Range("E51").Value = "  -0.53%  "
